$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 DecisionTreeClassifier(max_depth=10, min_samples_split=5,
                                        random_state=42))])'
$ws.Range("B2").Value = 0.6504778554778554
$ws.Range("C2").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': MinMaxScaler(), ''model__min_samples_split'': 5, ''model__min_samples_leaf'': 1, ''model__max_features'': None, ''model__max_depth'': 10, ''model__criterion'': ''gini'', ''model__class_weight'': None}'
$ws.Range("D2").Value = 0.7110540026094299
$ws.Range("E2").Value = 0.4522316394716394
$ws.Range("F2").Value = 0.5185185185185185
$ws.Range("G2").Value = 0.8634033027666478
$ws.Range("H2").Value = 0.539687619047619
$ws.Range("I2").Value = 0.6363636363636364
$ws.Range("J2").Value = 0.6444340425531915
$ws.Range("K2").Value = 0.4296666666666666
$ws.Range("L2").Value = 0.4375
$ws.Range("M2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range("N2").Value = '[0 1 1 1 1 1 0 0 1 0 0 1 0 0 0 0 1 1 0 0 1 1 0 0]'
$ws.Range("O2").Value = 42
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 DecisionTreeClassifier(class_weight=''balanced'',
                                        criterion=''entropy'', max_depth=15,
                                        max_features=''log2'',
                                        min_samples_split=7,
                                        random_state=42))])'
$ws.Range("B3").Value = 0.5987012987012987
$ws.Range("C3").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': RobustScaler(), ''model__min_samples_split'': 7, ''model__min_samples_leaf'': 1, ''model__max_features'': ''log2'', ''model__max_depth'': 15, ''model__criterion'': ''entropy'', ''model__class_weight'': ''balanced''}'
$ws.Range("D3").Value = 0.7044712316654341
$ws.Range("E3").Value = 0.4406009879009879
$ws.Range("F3").Value = 0.5185185185185185
$ws.Range("G3").Value = 0.8777631029897763
$ws.Range("H3").Value = 0.6179201587301587
$ws.Range("I3").Value = 0.6363636363636364
$ws.Range("J3").Value = 0.6313702127659574
$ws.Range("K3").Value = 0.3791333333333333
$ws.Range("L3").Value = 0.4375
$ws.Range("M3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]'
$ws.Range("N3").Value = '[0 1 1 0 0 1 1 0 1 0 0 0 0 0 1 1 1 0 0 1 1 0 0 1]'
$ws.Range("O3").Value = 69
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 DecisionTreeClassifier(class_weight=''balanced'',
                                        criterion=''entropy'', max_depth=3,
                                        max_features=''log2'', min_samples_leaf=7,
                                        min_samples_split=9,
                                        random_state=42))])'
$ws.Range("B4").Value = 0.6486596736596737
$ws.Range("C4").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': RobustScaler(), ''model__min_samples_split'': 9, ''model__min_samples_leaf'': 7, ''model__max_features'': ''log2'', ''model__max_depth'': 3, ''model__criterion'': ''entropy'', ''model__class_weight'': ''balanced''}'
$ws.Range("D4").Value = 0.7558565302625538
$ws.Range("E4").Value = 0.491752947052947
$ws.Range("F4").Value = 0.5714285714285714
$ws.Range("G4").Value = 0.8480025757893054
$ws.Range("H4").Value = 0.5230658730158729
$ws.Range("I4").Value = 0.8888888888888888
$ws.Range("J4").Value = 0.7155466666666667
$ws.Range("K4").Value = 0.5034399999999999
$ws.Range("L4").Value = 0.4210526315789473
$ws.Range("M4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range("N4").Value = '[0 0 0 0 0 1 0 1 0 0 1 0 1 1 1 0 0 0 0 0 1 0 1 1]'
$ws.Range("O4").Value = 23
$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 DecisionTreeClassifier(class_weight=''balanced'', max_depth=10,
                                        max_features=''sqrt'', min_samples_leaf=9,
                                        min_samples_split=7,
                                        random_state=42))])'
$ws.Range("B5").Value = 0.5770695970695972
$ws.Range("C5").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': StandardScaler(), ''model__min_samples_split'': 7, ''model__min_samples_leaf'': 9, ''model__max_features'': ''sqrt'', ''model__max_depth'': 10, ''model__criterion'': ''gini'', ''model__class_weight'': ''balanced''}'
$ws.Range("D5").Value = 0.7125511372454366
$ws.Range("E5").Value = 0.4530666156066156
$ws.Range("F5").Value = 0.5384615384615384
$ws.Range("G5").Value = 0.8607701385465948
$ws.Range("H5").Value = 0.5649090476190477
$ws.Range("I5").Value = 0.5833333333333334
$ws.Range("J5").Value = 0.6429551020408163
$ws.Range("K5").Value = 0.4134666666666667
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = '[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]'
$ws.Range("N5").Value = '[0 1 1 1 0 0 1 1 0 1 1 1 0 0 1 0 0 0 0 0 0 1 1 1]'
$ws.Range("O5").Value = 99
$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 DecisionTreeClassifier(criterion=''entropy'', max_depth=1,
                                        min_samples_leaf=9, min_samples_split=7,
                                        random_state=42))])'
$ws.Range("B6").Value = 0.6257420357420357
$ws.Range("C6").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': MinMaxScaler(), ''model__min_samples_split'': 7, ''model__min_samples_leaf'': 9, ''model__max_features'': None, ''model__max_depth'': 1, ''model__criterion'': ''entropy'', ''model__class_weight'': None}'
$ws.Range("D6").Value = 0.7479948781788192
$ws.Range("E6").Value = 0.5167622577422577
$ws.Range("F6").Value = 0.5925925925925926
$ws.Range("G6").Value = 0.8744301934939428
$ws.Range("H6").Value = 0.6382023809523809
$ws.Range("I6").Value = 0.5
$ws.Range("J6").Value = 0.668823076923077
$ws.Range("K6").Value = 0.4762666666666667
$ws.Range("L6").Value = 0.7272727272727273
$ws.Range("M6").Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]'
$ws.Range("N6").Value = '[1 1 1 1 1 1 0 0 0 0 1 0 1 1 1 0 0 1 0 1 1 1 1 1]'
$ws.Range("O6").Value = 89
